$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.1772000188674474
$ws.Range("J4").Value = 0.4772399240310436
$ws.Range("K4").Value = 0.7597635802022941
$ws.Range("L4").Value = 3.304213389742004
